$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.019.13"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "3.384.32"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.99"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.90"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.65"
$ws.Range("E9").Value = "  +2.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.122"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "3.962.89"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.86"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "3.390.49"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D17").Value = "61.117.01"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.11"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.61"
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.78"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.79"
$ws.Range("E22").Value = "  +3.57%  "
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "3.520.21"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("E27").Value = "  +6.99%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.14"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.37"
$ws.Range("E33").Value = "  -2.98%  "
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.95"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.74"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("D37").Value = "3.418.95"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("E39").Value = "  -2.58%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.17"
$ws.Range("E41").Value = "  -2.42%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.780"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.36"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.64"
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("D47").Value = "2.446.74"
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.97"
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.64"
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("E50").Value = "  +9.91%  "
$ws.Range("E51").Value = "  -2.09%  "
